# Automatic update of files.
# - "Förändrad" (column C) date serial for rows 2-8 bumps from 46078 to 46079.
# - Rows 6 and 7 swap their Beteckning/Datum/Area values (A, B, G columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump the "Förändrad" column (C) for rows 2 through 8 to the new date serial.
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 3).Value = 46079
}

# Swap the contents of rows 6 and 7 for columns A (Beteckning), B (Datum), and G (Area).
# Read values first (reading uses the Value() accessor form in this host),
# then write them after all reads are captured so row 6's write doesn't
# clobber what we still need to read from row 7.
$a6 = $ws.Cells.Item(6, 1).Value()
$b6 = $ws.Cells.Item(6, 2).Value()
$g6 = $ws.Cells.Item(6, 7).Value()

$a7 = $ws.Cells.Item(7, 1).Value()
$b7 = $ws.Cells.Item(7, 2).Value()
$g7 = $ws.Cells.Item(7, 7).Value()

$ws.Cells.Item(6, 1).Value = $a7
$ws.Cells.Item(6, 2).Value = $b7
$ws.Cells.Item(6, 7).Value = $g7

$ws.Cells.Item(7, 1).Value = $a6
$ws.Cells.Item(7, 2).Value = $b6
$ws.Cells.Item(7, 7).Value = $g6
